# Update the login credentials on the active sheet:
#   A2 (usuario): mauricio -> teste
#   B2 (senha):   3081995  -> teste   (now stored as text, not a number)
#   C2 (perfil):  admin    (unchanged value, kept as-is)
# Also move the active selection from C2 to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "teste"
$ws.Range("B2").Value = "teste"
$ws.Range("C2").Value = "admin"

$ws.Range("A3").Select()
